{"js": "// Update each multiplication-problem cell in the practice table to the\n// new operands/operator shown in the commit's target OOXML. Every\n// \"NNN\u00d7N=\" expression in the document is replaced by a new one (25\n// cells total); each old string is unique in the document, so a plain\n// matchCase search + full-match replace is safe and unambiguous.\nconst replacements = [\n  [\"512\u00d73=\", \"633\u00d74=\"],\n  [\"221\u00d72=\", \"325\u00d79=\"],\n  [\"206\u00d74=\", \"176\u00d77=\"],\n  [\"543\u00d74=\", \"345\u00d78=\"],\n  [\"775\u00d78=\", \"176\u00d76=\"],\n  [\"991\u00d75=\", \"854\u00d78=\"],\n  [\"730\u00d76=\", \"234\u00d75=\"],\n  [\"655\u00d74=\", \"452\u00d79=\"],\n  [\"333\u00d77=\", \"423\u00d77=\"],\n  [\"648\u00d78=\", \"669\u00d74=\"],\n  [\"257\u00d77=\", \"727\u00d76=\"],\n  [\"497\u00d78=\", \"188\u00d74=\"],\n  [\"604\u00d77=\", \"668\u00d76=\"],\n  [\"739\u00d73=\", \"735\u00d73=\"],\n  [\"436\u00d74=\", \"554\u00d72=\"],\n  [\"451\u00d78=\", \"713\u00d72=\"],\n  [\"986\u00d77=\", \"139\u00d72=\"],\n  [\"659\u00d78=\", \"190\u00d79=\"],\n  [\"530\u00d76=\", \"504\u00d74=\"],\n  [\"738\u00d73=\", \"866\u00d76=\"],\n  [\"346\u00d78=\", \"635\u00d79=\"],\n  [\"988\u00d76=\", \"873\u00d79=\"],\n  [\"721\u00d73=\", \"552\u00d77=\"],\n  [\"608\u00d76=\", \"499\u00d76=\"],\n  [\"254\u00d76=\", \"586\u00d77=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const found = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  found.load(\"items\");\n  await context.sync();\n\n  if (found.items.length === 0) {\n    throw new Error(`Text not found: ${oldText}`);\n  }\n\n  for (const range of found.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Update each multiplication-problem cell in the practice table to the\n# new operands/operator shown in the commit's target OOXML. Every\n# \"NNN\u00d7N=\" expression in the document is replaced by a new one (25\n# cells total); each old string is unique in the document, so a plain\n# Find/Replace (match whole string, no wildcards) for each pair is safe\n# and unambiguous.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"512\u00d73=\", \"633\u00d74=\"),\n    @(\"221\u00d72=\", \"325\u00d79=\"),\n    @(\"206\u00d74=\", \"176\u00d77=\"),\n    @(\"543\u00d74=\", \"345\u00d78=\"),\n    @(\"775\u00d78=\", \"176\u00d76=\"),\n    @(\"991\u00d75=\", \"854\u00d78=\"),\n    @(\"730\u00d76=\", \"234\u00d75=\"),\n    @(\"655\u00d74=\", \"452\u00d79=\"),\n    @(\"333\u00d77=\", \"423\u00d77=\"),\n    @(\"648\u00d78=\", \"669\u00d74=\"),\n    @(\"257\u00d77=\", \"727\u00d76=\"),\n    @(\"497\u00d78=\", \"188\u00d74=\"),\n    @(\"604\u00d77=\", \"668\u00d76=\"),\n    @(\"739\u00d73=\", \"735\u00d73=\"),\n    @(\"436\u00d74=\", \"554\u00d72=\"),\n    @(\"451\u00d78=\", \"713\u00d72=\"),\n    @(\"986\u00d77=\", \"139\u00d72=\"),\n    @(\"659\u00d78=\", \"190\u00d79=\"),\n    @(\"530\u00d76=\", \"504\u00d74=\"),\n    @(\"738\u00d73=\", \"866\u00d76=\"),\n    @(\"346\u00d78=\", \"635\u00d79=\"),\n    @(\"988\u00d76=\", \"873\u00d79=\"),\n    @(\"721\u00d73=\", \"552\u00d77=\"),\n    @(\"608\u00d76=\", \"499\u00d76=\"),\n    @(\"254\u00d76=\", \"586\u00d77=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Forward = $true\n    $find.Wrap = 0\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.Execute($find.Text, $find.MatchCase, $find.MatchWholeWord, $find.MatchWildcards, $false, $false, $find.Forward, $find.Wrap, $false, $find.Replacement.Text, 2)\n}\n"}
